$d = $word.ActiveDocument

# Locate the "SMARTREWARDS FAQ'S" heading paragraph and the final (trailing empty)
# paragraph in the document; delete everything in between (the whole FAQ block),
# leaving the heading paragraph and the trailing empty paragraph intact.

$headingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "SMARTREWARDS FAQ*") {
        $headingPara = $p
    }
}

$lastPara = $d.Paragraphs.Last

$startPos = $headingPara.Range.End
$endPos = $lastPara.Range.Start

$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()
